$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

# --- Neg_Change sheet: update rows 2-14 in place (row count unchanged) ---
$ws1.Cells.Item(2, 1).Value = "DRREDDY"
$ws1.Cells.Item(2, 2).Value = 1250
$ws1.Cells.Item(2, 3).Value = 1258.8
$ws1.Cells.Item(2, 4).Value = 1245.5
$ws1.Cells.Item(2, 5).Value = 1255
$ws1.Cells.Item(2, 6).Value = 540944
$ws1.Cells.Item(2, 7).Value = 1184386
$ws1.Cells.Item(2, 8).Value = -0.5432705216035988
$ws1.Cells.Item(2, 9).Value = "DRREDDY"
$ws1.Cells.Item(3, 1).Value = "NTPC"
$ws1.Cells.Item(3, 2).Value = 349.9
$ws1.Cells.Item(3, 3).Value = 354.1
$ws1.Cells.Item(3, 4).Value = 346.15
$ws1.Cells.Item(3, 5).Value = 350.5
$ws1.Cells.Item(3, 6).Value = 8641743
$ws1.Cells.Item(3, 7).Value = 18248271
$ws1.Cells.Item(3, 8).Value = -0.5264349701952585
$ws1.Cells.Item(3, 9).Value = "NTPC"
$ws1.Cells.Item(4, 1).Value = "SOLARINDS"
$ws1.Cells.Item(4, 2).Value = 12749
$ws1.Cells.Item(4, 3).Value = 12870
$ws1.Cells.Item(4, 4).Value = 12585
$ws1.Cells.Item(4, 5).Value = 12850
$ws1.Cells.Item(4, 6).Value = 95209
$ws1.Cells.Item(4, 7).Value = 231189
$ws1.Cells.Item(4, 8).Value = -0.588176773116368
$ws1.Cells.Item(4, 9).Value = "SOLARINDS"
$ws1.Cells.Item(5, 1).Value = "LICI"
$ws1.Cells.Item(5, 2).Value = 848
$ws1.Cells.Item(5, 3).Value = 858.6
$ws1.Cells.Item(5, 4).Value = 848
$ws1.Cells.Item(5, 5).Value = 851.95
$ws1.Cells.Item(5, 6).Value = 969130
$ws1.Cells.Item(5, 7).Value = 2084331
$ws1.Cells.Item(5, 8).Value = -0.5350402599203293
$ws1.Cells.Item(5, 9).Value = "LICI"
$ws1.Cells.Item(6, 1).Value = "HAL"
$ws1.Cells.Item(6, 2).Value = 4528.3
$ws1.Cells.Item(6, 3).Value = 4547.7
$ws1.Cells.Item(6, 4).Value = 4478.6
$ws1.Cells.Item(6, 5).Value = 4505.1
$ws1.Cells.Item(6, 6).Value = 530874
$ws1.Cells.Item(6, 7).Value = 1256096
$ws1.Cells.Item(6, 8).Value = -0.5773619213818052
$ws1.Cells.Item(6, 9).Value = "HAL"
$ws1.Cells.Item(7, 1).Value = "MAZDOCK"
$ws1.Cells.Item(7, 2).Value = 2520
$ws1.Cells.Item(7, 3).Value = 2520.4
$ws1.Cells.Item(7, 4).Value = 2473.8
$ws1.Cells.Item(7, 5).Value = 2493.9
$ws1.Cells.Item(7, 6).Value = 582637
$ws1.Cells.Item(7, 7).Value = 1160745
$ws1.Cells.Item(7, 8).Value = -0.4980490977777203
$ws1.Cells.Item(7, 9).Value = "MAZDOCK"
$ws1.Cells.Item(8, 1).Value = "AMBUJACEM"
$ws1.Cells.Item(8, 2).Value = 572
$ws1.Cells.Item(8, 3).Value = 572
$ws1.Cells.Item(8, 4).Value = 561.35
$ws1.Cells.Item(8, 5).Value = 563
$ws1.Cells.Item(8, 6).Value = 1015560
$ws1.Cells.Item(8, 7).Value = 2369536
$ws1.Cells.Item(8, 8).Value = -0.5714097612359551
$ws1.Cells.Item(8, 9).Value = "AMBUJACEM"
$ws1.Cells.Item(9, 1).Value = "TIINDIA"
$ws1.Cells.Item(9, 2).Value = 2530
$ws1.Cells.Item(9, 3).Value = 2560
$ws1.Cells.Item(9, 4).Value = 2517.1
$ws1.Cells.Item(9, 5).Value = 2549.9
$ws1.Cells.Item(9, 6).Value = 240876
$ws1.Cells.Item(9, 7).Value = 501042
$ws1.Cells.Item(9, 8).Value = -0.5192498832433209
$ws1.Cells.Item(9, 9).Value = "TIINDIA"
$ws1.Cells.Item(10, 1).Value = "PRESTIGE"
$ws1.Cells.Item(10, 2).Value = 1675.4
$ws1.Cells.Item(10, 3).Value = 1676.9
$ws1.Cells.Item(10, 4).Value = 1648
$ws1.Cells.Item(10, 5).Value = 1657.3
$ws1.Cells.Item(10, 6).Value = 448673
$ws1.Cells.Item(10, 7).Value = 906033
$ws1.Cells.Item(10, 8).Value = -0.5047939754953737
$ws1.Cells.Item(10, 9).Value = "PRESTIGE"
$ws1.Cells.Item(11, 1).Value = "PAGEIND"
$ws1.Cells.Item(11, 2).Value = 35685
$ws1.Cells.Item(11, 3).Value = 35685
$ws1.Cells.Item(11, 4).Value = 34990
$ws1.Cells.Item(11, 5).Value = 35300
$ws1.Cells.Item(11, 6).Value = 11501
$ws1.Cells.Item(11, 7).Value = 27388
$ws1.Cells.Item(11, 8).Value = -0.5800715641886958
$ws1.Cells.Item(11, 9).Value = "PAGEIND"
$ws1.Cells.Item(12, 1).Value = "UNOMINDA"
$ws1.Cells.Item(12, 2).Value = 1342.6
$ws1.Cells.Item(12, 3).Value = 1343.1
$ws1.Cells.Item(12, 4).Value = 1307.9
$ws1.Cells.Item(12, 5).Value = 1320
$ws1.Cells.Item(12, 6).Value = 396153
$ws1.Cells.Item(12, 7).Value = 918558
$ws1.Cells.Item(12, 8).Value = -0.5687229331190845
$ws1.Cells.Item(12, 9).Value = "UNOMINDA"
$ws1.Cells.Item(13, 1).Value = "BLUESTARCO"
$ws1.Cells.Item(13, 2).Value = 1850
$ws1.Cells.Item(13, 3).Value = 1854.5
$ws1.Cells.Item(13, 4).Value = 1806
$ws1.Cells.Item(13, 5).Value = 1823.4
$ws1.Cells.Item(13, 6).Value = 266141
$ws1.Cells.Item(13, 7).Value = 602390
$ws1.Cells.Item(13, 8).Value = -0.5581915370441076
$ws1.Cells.Item(13, 9).Value = "BLUESTARCO"
$ws1.Cells.Item(14, 1).Value = "GMRAIRPORT"
$ws1.Cells.Item(14, 2).Value = 105.74
$ws1.Cells.Item(14, 3).Value = 105.77
$ws1.Cells.Item(14, 4).Value = 103.37
$ws1.Cells.Item(14, 5).Value = 104.25
$ws1.Cells.Item(14, 6).Value = 5695438
$ws1.Cells.Item(14, 7).Value = 11284687
$ws1.Cells.Item(14, 8).Value = -0.4952949957761345
$ws1.Cells.Item(14, 9).Value = "GMRAIRPORT"

# --- Pos_Change sheet: update rows 2-15 in place ---
$ws2.Cells.Item(2, 1).Value = "JIOFIN"
$ws2.Cells.Item(2, 2).Value = 300.1
$ws2.Cells.Item(2, 3).Value = 302.3
$ws2.Cells.Item(2, 4).Value = 293.5
$ws2.Cells.Item(2, 5).Value = 298.55
$ws2.Cells.Item(2, 6).Value = 10930694
$ws2.Cells.Item(2, 7).Value = 6971528
$ws2.Cells.Item(2, 8).Value = 0.5679050561082162
$ws2.Cells.Item(2, 9).Value = "JIOFIN"
$ws2.Cells.Item(3, 1).Value = "ETERNAL"
$ws2.Cells.Item(3, 2).Value = 281.8
$ws2.Cells.Item(3, 3).Value = 283.55
$ws2.Cells.Item(3, 4).Value = 276.25
$ws2.Cells.Item(3, 5).Value = 279
$ws2.Cells.Item(3, 6).Value = 29577754
$ws2.Cells.Item(3, 7).Value = 19328131
$ws2.Cells.Item(3, 8).Value = 0.5302956090270704
$ws2.Cells.Item(3, 9).Value = "ETERNAL"
$ws2.Cells.Item(4, 1).Value = "ITC"
$ws2.Cells.Item(4, 2).Value = 350.25
$ws2.Cells.Item(4, 3).Value = 350.7
$ws2.Cells.Item(4, 4).Value = 337.75
$ws2.Cells.Item(4, 5).Value = 343.25
$ws2.Cells.Item(4, 6).Value = 72582738
$ws2.Cells.Item(4, 7).Value = 48358124
$ws2.Cells.Item(4, 8).Value = 0.5009419720252175
$ws2.Cells.Item(4, 9).Value = "ITC"
$ws2.Cells.Item(5, 1).Value = "SHREECEM"
$ws2.Cells.Item(5, 2).Value = 27670
$ws2.Cells.Item(5, 3).Value = 28090
$ws2.Cells.Item(5, 4).Value = 27545
$ws2.Cells.Item(5, 5).Value = 27615
$ws2.Cells.Item(5, 6).Value = 35574
$ws2.Cells.Item(5, 7).Value = 23089
$ws2.Cells.Item(5, 8).Value = 0.5407336827060505
$ws2.Cells.Item(5, 9).Value = "SHREECEM"
$ws2.Cells.Item(6, 1).Value = "INDHOTEL"
$ws2.Cells.Item(6, 2).Value = 746
$ws2.Cells.Item(6, 3).Value = 746
$ws2.Cells.Item(6, 4).Value = 720.3
$ws2.Cells.Item(6, 5).Value = 725.95
$ws2.Cells.Item(6, 6).Value = 2525206
$ws2.Cells.Item(6, 7).Value = 1693043
$ws2.Cells.Item(6, 8).Value = 0.4915191167619488
$ws2.Cells.Item(6, 9).Value = "INDHOTEL"
$ws2.Cells.Item(7, 1).Value = "ICICIPRULI"
$ws2.Cells.Item(7, 2).Value = 685.05
$ws2.Cells.Item(7, 3).Value = 697.5
$ws2.Cells.Item(7, 4).Value = 685.05
$ws2.Cells.Item(7, 5).Value = 690
$ws2.Cells.Item(7, 6).Value = 920286
$ws2.Cells.Item(7, 7).Value = 579420
$ws2.Cells.Item(7, 8).Value = 0.5882882882882883
$ws2.Cells.Item(7, 9).Value = "ICICIPRULI"
$ws2.Cells.Item(8, 1).Value = "BDL"
$ws2.Cells.Item(8, 2).Value = 1540
$ws2.Cells.Item(8, 3).Value = 1564.7
$ws2.Cells.Item(8, 4).Value = 1524
$ws2.Cells.Item(8, 5).Value = 1542
$ws2.Cells.Item(8, 6).Value = 4768472
$ws2.Cells.Item(8, 7).Value = 3329765
$ws2.Cells.Item(8, 8).Value = 0.4320746358977285
$ws2.Cells.Item(8, 9).Value = "BDL"
$ws2.Cells.Item(9, 1).Value = "IRB"
$ws2.Cells.Item(9, 2).Value = 42.25
$ws2.Cells.Item(9, 3).Value = 42.45
$ws2.Cells.Item(9, 4).Value = 41.79
$ws2.Cells.Item(9, 5).Value = 41.99
$ws2.Cells.Item(9, 6).Value = 7654980
$ws2.Cells.Item(9, 7).Value = 5006189
$ws2.Cells.Item(9, 8).Value = 0.5291032759650105
$ws2.Cells.Item(9, 9).Value = "IRB"
$ws2.Cells.Item(10, 1).Value = "SONACOMS"
$ws2.Cells.Item(10, 2).Value = 481.4
$ws2.Cells.Item(10, 3).Value = 481.5
$ws2.Cells.Item(10, 4).Value = 469.1
$ws2.Cells.Item(10, 5).Value = 476.2
$ws2.Cells.Item(10, 6).Value = 1057436
$ws2.Cells.Item(10, 7).Value = 665079
$ws2.Cells.Item(10, 8).Value = 0.5899404431653984
$ws2.Cells.Item(10, 9).Value = "SONACOMS"
$ws2.Cells.Item(11, 1).Value = "OFSS"
$ws2.Cells.Item(11, 2).Value = 7641
$ws2.Cells.Item(11, 3).Value = 7701
$ws2.Cells.Item(11, 4).Value = 7574.5
$ws2.Cells.Item(11, 5).Value = 7582.5
$ws2.Cells.Item(11, 6).Value = 160627
$ws2.Cells.Item(11, 7).Value = 104098
$ws2.Cells.Item(11, 8).Value = 0.5430363695748237
$ws2.Cells.Item(11, 9).Value = "OFSS"
$ws2.Cells.Item(12, 1).Value = "HFCL"
$ws2.Cells.Item(12, 2).Value = 67.45
$ws2.Cells.Item(12, 3).Value = 70.09999999999999
$ws2.Cells.Item(12, 4).Value = 66.34999999999999
$ws2.Cells.Item(12, 5).Value = 69.75
$ws2.Cells.Item(12, 6).Value = 26747246
$ws2.Cells.Item(12, 7).Value = 18916545
$ws2.Cells.Item(12, 8).Value = 0.4139604245912771
$ws2.Cells.Item(12, 9).Value = "HFCL"
$ws2.Cells.Item(13, 1).Value = "MCX"
$ws2.Cells.Item(13, 2).Value = 2212
$ws2.Cells.Item(13, 3).Value = 2270
$ws2.Cells.Item(13, 4).Value = 2190
$ws2.Cells.Item(13, 5).Value = 2249
$ws2.Cells.Item(13, 6).Value = 2414137
$ws2.Cells.Item(13, 7).Value = 1638124
$ws2.Cells.Item(13, 8).Value = 0.4737205486275764
$ws2.Cells.Item(13, 9).Value = "MCX"
$ws2.Cells.Item(14, 1).Value = "ABFRL"
$ws2.Cells.Item(14, 2).Value = 77.04000000000001
$ws2.Cells.Item(14, 3).Value = 77.81
$ws2.Cells.Item(14, 4).Value = 75.90000000000001
$ws2.Cells.Item(14, 5).Value = 77.27
$ws2.Cells.Item(14, 6).Value = 2364806
$ws2.Cells.Item(14, 7).Value = 1677318
$ws2.Cells.Item(14, 8).Value = 0.4098733811954561
$ws2.Cells.Item(14, 9).Value = "ABFRL"
$ws2.Cells.Item(15, 1).Value = "KFINTECH"
$ws2.Cells.Item(15, 2).Value = 1085.2
$ws2.Cells.Item(15, 3).Value = 1085.5
$ws2.Cells.Item(15, 4).Value = 1065.5
$ws2.Cells.Item(15, 5).Value = 1076.8
$ws2.Cells.Item(15, 6).Value = 1402596
$ws2.Cells.Item(15, 7).Value = 945849
$ws2.Cells.Item(15, 8).Value = 0.4828963185455606
$ws2.Cells.Item(15, 9).Value = "KFINTECH"

# --- Pos_Change sheet: remove now-obsolete rows 16-22 ---
$ws2.Rows("16:22").Delete()
